$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Year")

# --- Rows 1-5: re-arrange the A/B "parameter" block, add bold style to A1:A5 and D1:D5 ---
# Row 1: A1 becomes "Obj" (sum of results), B1 becomes the grand-total formula
$ws.Range("A1").Value = "Obj"
$ws.Range("B1").Formula = "=SUM(D9:D20)+SUM(K24:K35)+SUM(G39:G50)"

# Row 2: A2/B2 become empty (the "eH" pair moves down to row 4)
$ws.Range("A2").ClearContents()
$ws.Range("B2").ClearContents()

# Row 3: now holds what used to be in row 1 ("dP")
$ws.Range("A3").Value = "dP"
$ws.Range("B3").Value = 5

# Row 4: now holds what used to be in row 2 ("eH")
$ws.Range("A4").Value = "eH"
$ws.Range("B4").Value = 70

# Row 5: now holds what used to be in row 3 ("H")
$ws.Range("A5").Value = "H"
$ws.Range("B5").Value = 720

# Bold styling for the whole A1:A5 / D1:D5 label block
$ws.Range("A1:A5").Font.Bold = $true
$ws.Range("D1:D5").Font.Bold = $true

# --- Rows 9-20: update the $B$1*$B$2 references to $B$3*$B$4, and add new turbine data ---
$ws.Range("B15").Value = 13
$ws.Range("B17").Value = 7

$ws.Range("D9").Formula = "=(B9+C9)*`$B`$3*`$B`$4"
$ws.Range("D10").Formula = "=(B10+C10)*`$B`$3*`$B`$4"
$ws.Range("D11").Formula = "=(B11+C11)*`$B`$3*`$B`$4"
$ws.Range("D12").Formula = "=(B12+C12)*`$B`$3*`$B`$4"
$ws.Range("D13").Formula = "=(B13+C13)*`$B`$3*`$B`$4"
$ws.Range("D14").Formula = "=(B14+C14)*`$B`$3*`$B`$4"
$ws.Range("D15").Formula = "=(B15+C15)*`$B`$3*`$B`$4"
$ws.Range("D16").Formula = "=(B16+C16)*`$B`$3*`$B`$4"
$ws.Range("D17").Formula = "=(B17+C17)*`$B`$3*`$B`$4"
$ws.Range("D18").Formula = "=(B18+C18)*`$B`$3*`$B`$4"
$ws.Range("D19").Formula = "=(B19+C19)*`$B`$3*`$B`$4"
$ws.Range("D20").Formula = "=(B20+C20)*`$B`$3*`$B`$4"

# --- Rows 24-35: update the $B$2*$B$3 references to $B$4*$B$5 for the I and J columns ---
$ws.Range("I24").Formula = "=SUM(B`$24:B24,E`$24:E24)*`$B`$4*`$B`$5"
$ws.Range("J24").Formula = "=SUM(C`$24:C24,F`$24:F24)*`$B`$4*`$B`$5"
$ws.Range("I25").Formula = "=SUM(B`$24:B25,E`$24:E25)*`$B`$4*`$B`$5"
$ws.Range("J25").Formula = "=SUM(C`$24:C25,F`$24:F25)*`$B`$4*`$B`$5"
$ws.Range("I26").Formula = "=SUM(B`$24:B26,E`$24:E26)*`$B`$4*`$B`$5"
$ws.Range("J26").Formula = "=SUM(C`$24:C26,F`$24:F26)*`$B`$4*`$B`$5"
$ws.Range("I27").Formula = "=SUM(B`$24:B27,E`$24:E27)*`$B`$4*`$B`$5"
$ws.Range("J27").Formula = "=SUM(C`$24:C27,F`$24:F27)*`$B`$4*`$B`$5"
$ws.Range("I28").Formula = "=SUM(B`$24:B28,E`$24:E28)*`$B`$4*`$B`$5"
$ws.Range("J28").Formula = "=SUM(C`$24:C28,F`$24:F28)*`$B`$4*`$B`$5"
$ws.Range("I29").Formula = "=SUM(B`$24:B29,E`$24:E29)*`$B`$4*`$B`$5"
$ws.Range("J29").Formula = "=SUM(C`$24:C29,F`$24:F29)*`$B`$4*`$B`$5"
$ws.Range("I30").Formula = "=SUM(B`$24:B30,E`$24:E30)*`$B`$4*`$B`$5"
$ws.Range("J30").Formula = "=SUM(C`$24:C30,F`$24:F30)*`$B`$4*`$B`$5"
$ws.Range("I31").Formula = "=SUM(B`$24:B31,E`$24:E31)*`$B`$4*`$B`$5"
$ws.Range("J31").Formula = "=SUM(C`$24:C31,F`$24:F31)*`$B`$4*`$B`$5"
$ws.Range("I32").Formula = "=SUM(B`$24:B32,E`$24:E32)*`$B`$4*`$B`$5"
$ws.Range("J32").Formula = "=SUM(C`$24:C32,F`$24:F32)*`$B`$4*`$B`$5"
$ws.Range("I33").Formula = "=SUM(B`$24:B33,E`$24:E33)*`$B`$4*`$B`$5"
$ws.Range("J33").Formula = "=SUM(C`$24:C33,F`$24:F33)*`$B`$4*`$B`$5"
$ws.Range("I34").Formula = "=SUM(B`$24:B34,E`$24:E34)*`$B`$4*`$B`$5"
$ws.Range("J34").Formula = "=SUM(C`$24:C34,F`$24:F34)*`$B`$4*`$B`$5"
$ws.Range("I35").Formula = "=SUM(B`$24:B35,E`$24:E35)*`$B`$4*`$B`$5"
$ws.Range("J35").Formula = "=SUM(C`$24:C35,F`$24:F35)*`$B`$4*`$B`$5"

# --- Rows 39-50: re-enter D39:D50 as one shared-formula block, add a new G39:G50 column ---
$ws.Range("D39:D50").Formula = "=B39*`$E`$1+C39"
$ws.Range("G39:G50").Formula = "=B39*`$E`$3"

# --- Conditional formatting on E39:E50: swap the two rules so "less than or equal" is evaluated first ---
$cf = $ws.Range("E39:E50").FormatConditions
$cf.Delete()
$ruleLE = $ws.Range("E39:E50").FormatConditions.Add(8, 7, "=`$D39")
$ruleLE.Interior.Color = 13011546
$ruleLE.Font.Color = 6299648
$ruleGT = $ws.Range("E39:E50").FormatConditions.Add(8, 5, "=`$D39")
$ruleGT.Interior.Color = 13551615
$ruleGT.Font.Color = 192

$ws.Range("K25").Select()
